# Apply "add nowcasts for 2025q4" changes to the single worksheet.
# This replaces the header row + first 6 data rows of the nowcasts table
# with the new 2025Q4 figures, and nudges a handful of column widths.
# Rows 8-11 (the tail of the old table) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -----------------------------------------------
$headerValues = @("Row","Prognose","surveys","production","orders","turnover","financial","labor market","prices","national accounts","Revision")
for ($col = 1; $col -le $headerValues.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $headerValues[$col - 1]
}

# --- Row labels (column A, rows 2-7) -----------------------------------
# Make sure these stay plain text (not auto-converted to date serials)
# by forcing a text number format before assignment, then restoring the
# default "Normal" style so the cell matches the rest of the sheet.
$labelRange = $ws.Range("A2:A7")
$labelRange.NumberFormat = "@"

$rowLabels = @("2025-09-30","2025-10-15","2025-10-30","2025-11-15","2025-11-30","2025-12-15")
for ($i = 0; $i -lt $rowLabels.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $rowLabels[$i]
}

$labelRange.Style = "Normal"

# --- Numeric data, rows 2-7, columns B-K --------------------------------
$data = @{
    2 = @(0.1870989194238093, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    3 = @(0.26142233486996241, 0, 0.032614750438115861, -0.00072713620633531979, -0.013973723395611778, 0.015236237868631802, -0.0016667751737385033, 0.010709418967483071, 0, 0.032130642947607979)
    4 = @(0.59316401630507554, 0.12906820827666571, 0, 0.00036253763538482051, 0.00015911557385103005, 0, -0.0067252346047485066, 0.18630175941768695, 0.021246882872129395, 0.0013284122641437079)
    5 = @(0.47950293674006661, 0, -0.047544948711517931, 0.0090004095035760741, -0.079461724000078782, -0.013995479247274206, -0.0045854681692314679, 0.01825037813366669, 0, 0.0046757529258506514)
    6 = @(0.12895657549026832, -0.27604803347619739, 0, 0.0020961709841476033, -0.004102025004714991, 0, -0.0077309319767963893, -0.066170965371753321, 0, 0.0014094235955162127)
    7 = @(0.14425854688922357, 0, 0.16852345041549277, -0.0060927282445517005, -0.12715867628718663, 0.0052931668063104153, 0, 0, 0, -0.0252632412911096)
}

foreach ($rowNum in $data.Keys) {
    $vals = $data[$rowNum]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($rowNum, $c + 2).Value = $vals[$c]
    }
}

# --- Column width tweaks -------------------------------------------------
# Target (OOXML "width" character units): C=13.25, D=14.25, G=14.25,
# H=15.25, I=14.25, J=15.05, K=14.65 ; columns A, B, E, F stay as-is.
$ws.Columns.Item(3).ColumnWidth = 12.333333333333334
$ws.Columns.Item(4).ColumnWidth = 13.333333333333334
$ws.Columns.Item(7).ColumnWidth = 13.333333333333334
$ws.Columns.Item(8).ColumnWidth = 14.333333333333334
$ws.Columns.Item(9).ColumnWidth = 13.333333333333334
$ws.Columns.Item(10).ColumnWidth = 14.166666666666666
$ws.Columns.Item(11).ColumnWidth = 13.833333333333334
